$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 535.4848691260194
$ws.Range("C2").Value = 378.1733609837286
$ws.Range("D2").Value = 324.3964720222395
$ws.Range("E2").Value = 294.2179878809634
$ws.Range("B3").Value = 631.3495691272635
$ws.Range("C3").Value = 447.9257722692942
$ws.Range("D3").Value = 382.9822906716718
$ws.Range("E3").Value = 350.602469332094
$ws.Range("B4").Value = 598.8871305916017
$ws.Range("C4").Value = 424.2749863900274
$ws.Range("D4").Value = 361.8241150246059
$ws.Range("E4").Value = 330.8530760416022
$ws.Range("B5").Value = 413.9612067934758
$ws.Range("C5").Value = 294.3585856825255
$ws.Range("D5").Value = 248.7800367579859
$ws.Range("E5").Value = 230.6370800733493
$ws.Range("B6").Value = 369.0734575106675
$ws.Range("C6").Value = 260.8286789352674
$ws.Range("D6").Value = 221.7938529479062
$ws.Range("E6").Value = 202.9059229700368
$ws.Range("B7").Value = 38.33206710804144
$ws.Range("C7").Value = 27.0678373726868
$ws.Range("D7").Value = 23.24229472104152
$ws.Range("E7").Value = 21.20280047310334
$ws.Range("B8").Value = 2113.850600438589
$ws.Range("C8").Value = 1495.657357800417
$ws.Range("D8").Value = 1295.065367782164
$ws.Range("E8").Value = 1168.291416864006
$ws.Range("B9").Value = 530.3805112132241
$ws.Range("C9").Value = 376.7627727885709
$ws.Range("D9").Value = 322.1087199759793
$ws.Range("E9").Value = 295.2365459766266
$ws.Range("B10").Value = 196.0768590149285
$ws.Range("C10").Value = 135.5174981927556
$ws.Range("D10").Value = 116.9382027048995
$ws.Range("E10").Value = 105.3879872271629
$ws.Range("B11").Value = 34.98568552358677
$ws.Range("C11").Value = 22.87120849966349
$ws.Range("D11").Value = 19.47136745963188
$ws.Range("E11").Value = 18.93101564496779
$ws.Range("B12").Value = 72.49519911682752
$ws.Range("C12").Value = 51.03096720047353
$ws.Range("D12").Value = 42.92710729652431
$ws.Range("E12").Value = 38.08525385411853
$ws.Range("B13").Value = 109.3799209453936
$ws.Range("C13").Value = 74.84677731521124
$ws.Range("D13").Value = 65.12455326262715
$ws.Range("E13").Value = 59.37250690201907
